# Update the two-digit multiplication problems in the single table.
# Several "old" values are duplicated across cells (e.g. "34×95=" appears
# twice but must become two different new values), so plain document-wide
# Find/Replace is ambiguous. Instead we address each cell directly by its
# (row, column) position in the table, which uniquely identifies every
# occurrence regardless of duplicate text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "83×36="; New = "77×97=" },
    @{ Row = 1;  Col = 2; Old = "39×95="; New = "24×57=" },
    @{ Row = 1;  Col = 3; Old = "60×51="; New = "76×43=" },
    @{ Row = 1;  Col = 4; Old = "83×23="; New = "82×97=" },
    @{ Row = 1;  Col = 5; Old = "94×31="; New = "56×70=" },

    @{ Row = 5;  Col = 1; Old = "29×40="; New = "67×56=" },
    @{ Row = 5;  Col = 2; Old = "45×44="; New = "38×95=" },
    @{ Row = 5;  Col = 3; Old = "34×95="; New = "28×97=" },
    @{ Row = 5;  Col = 4; Old = "22×70="; New = "33×40=" },
    @{ Row = 5;  Col = 5; Old = "14×54="; New = "98×80=" },

    @{ Row = 10; Col = 1; Old = "77×27="; New = "84×61=" },
    @{ Row = 10; Col = 2; Old = "25×72="; New = "75×99=" },
    @{ Row = 10; Col = 3; Old = "26×50="; New = "81×44=" },
    @{ Row = 10; Col = 4; Old = "88×45="; New = "27×87=" },
    @{ Row = 10; Col = 5; Old = "54×54="; New = "83×44=" },

    @{ Row = 15; Col = 1; Old = "20×50="; New = "98×96=" },
    @{ Row = 15; Col = 2; Old = "13×87="; New = "58×45=" },
    @{ Row = 15; Col = 3; Old = "34×95="; New = "77×35=" },
    @{ Row = 15; Col = 4; Old = "76×22="; New = "62×99=" },
    @{ Row = 15; Col = 5; Old = "63×42="; New = "88×16=" },

    @{ Row = 20; Col = 1; Old = "46×49="; New = "50×18=" },
    @{ Row = 20; Col = 2; Old = "78×19="; New = "73×38=" },
    @{ Row = 20; Col = 3; Old = "58×32="; New = "31×22=" },
    @{ Row = 20; Col = 4; Old = "63×95="; New = "49×33=" },
    @{ Row = 20; Col = 5; Old = "56×41="; New = "26×13=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Replace just the single match within this cell's range (wdReplaceOne
    # = 1, not wdReplaceAll = 2). Several "old" expressions are duplicated
    # across multiple cells but map to different "new" values, and a
    # document-wide ReplaceAll would clobber every occurrence at once, so
    # each cell must be patched individually with a single replacement.
    $rng.Find.Execute($u.Old, $true, $false, $false, $false, $false, `
                       $true, 0, $false, $u.New, 1)
}
